{"js": "// The document is a title line (\"YYYY-MM-DD Weekday\") followed by a\n// table of two-digit \u00f7 one-digit division problems. Every non-empty\n// paragraph's text changes to a new value, in document order. Using\n// positional (paragraph-order) replacement avoids any ambiguity from\n// values that repeat between the \"before\" and \"after\" sets.\nconst replacements = [\n  \"2025-09-13 Saturday\",\n  \"74\u00f76=12, 2\",\n  \"72\u00f77=10, 2\",\n  \"73\u00f78=9, 1\",\n  \"81\u00f73=27, 0\",\n  \"80\u00f73=26, 2\",\n  \"93\u00f76=15, 3\",\n  \"40\u00f75=8, 0\",\n  \"21\u00f79=2, 3\",\n  \"23\u00f75=4, 3\",\n  \"35\u00f78=4, 3\",\n  \"89\u00f76=14, 5\",\n  \"57\u00f79=6, 3\",\n  \"23\u00f72=11, 1\",\n  \"34\u00f77=4, 6\",\n  \"51\u00f77=7, 2\",\n  \"60\u00f79=6, 6\",\n  \"23\u00f79=2, 5\",\n  \"92\u00f73=30, 2\",\n  \"84\u00f72=42, 0\",\n  \"33\u00f75=6, 3\",\n  \"33\u00f75=6, 3\",\n  \"10\u00f74=2, 2\",\n  \"79\u00f73=26, 1\",\n  \"36\u00f76=6, 0\",\n  \"98\u00f78=12, 2\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idx = 0;\nfor (const paragraph of paragraphs.items) {\n  if (idx >= replacements.length) break;\n  const text = paragraph.text;\n  // skip paragraphs with no visible text (the blank spacer rows in the table)\n  if (!text || text.trim().length === 0) continue;\n  paragraph.insertText(replacements[idx], Word.InsertLocation.replace);\n  idx++;\n}\n\nawait context.sync();\n", "ps1": "# The document is a title line (\"YYYY-MM-DD Weekday\") followed by a\n# table of two-digit \u00f7 one-digit division problems. Every non-empty\n# paragraph's text changes to a new value, in document order. Walking\n# the Paragraphs collection positionally (instead of searching by the\n# old text) avoids any ambiguity from values that repeat between the\n# \"before\" and \"after\" sets.\n$replacements = @(\n    \"2025-09-13 Saturday\",\n    \"74\u00f76=12, 2\",\n    \"72\u00f77=10, 2\",\n    \"73\u00f78=9, 1\",\n    \"81\u00f73=27, 0\",\n    \"80\u00f73=26, 2\",\n    \"93\u00f76=15, 3\",\n    \"40\u00f75=8, 0\",\n    \"21\u00f79=2, 3\",\n    \"23\u00f75=4, 3\",\n    \"35\u00f78=4, 3\",\n    \"89\u00f76=14, 5\",\n    \"57\u00f79=6, 3\",\n    \"23\u00f72=11, 1\",\n    \"34\u00f77=4, 6\",\n    \"51\u00f77=7, 2\",\n    \"60\u00f79=6, 6\",\n    \"23\u00f79=2, 5\",\n    \"92\u00f73=30, 2\",\n    \"84\u00f72=42, 0\",\n    \"33\u00f75=6, 3\",\n    \"33\u00f75=6, 3\",\n    \"10\u00f74=2, 2\",\n    \"79\u00f73=26, 1\",\n    \"36\u00f76=6, 0\",\n    \"98\u00f78=12, 2\"\n)\n\n$d = $word.ActiveDocument\n$idx = 0\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($idx -ge $replacements.Count) { break }\n    $para = $d.Paragraphs($i)\n    $text = $para.Range.Text -replace \"[\\r\\a\\n]\", \"\"\n    if ($text.Trim().Length -gt 0) {\n        $para.Range.Text = $replacements[$idx]\n        $idx++\n    }\n}\n"}
